$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 313; this shifts the existing rows
# 313-320 down to 315-322, carrying their content/formatting with them.
$ws.Rows.Item(313).Insert()
$ws.Rows.Item(313).Insert()

# Populate the first new row (313) - Lapins / Especial
$ws.Cells.Item(313, 1).Value  = 9
$ws.Cells.Item(313, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(313, 3).Value  = "Metropolitana"
$ws.Cells.Item(313, 4).Value  = 44595
$ws.Cells.Item(313, 5).Value  = 13
$ws.Cells.Item(313, 6).Value  = "Fruta"
$ws.Cells.Item(313, 7).Value  = 100103
$ws.Cells.Item(313, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(313, 9).Value  = 100103001
$ws.Cells.Item(313, 10).Value = "Cereza"
$ws.Cells.Item(313, 11).Value = "Lapins"
$ws.Cells.Item(313, 12).Value = "Especial"
$ws.Cells.Item(313, 13).Value = 200
$ws.Cells.Item(313, 14).Value = 12000
$ws.Cells.Item(313, 15).Value = 12000
$ws.Cells.Item(313, 16).Value = 12000
$ws.Cells.Item(313, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(313, 18).Value = "Región del Maule"
$ws.Cells.Item(313, 19).Value = 1200
$ws.Cells.Item(313, 20).Value = 10

# Populate the second new row (314) - Lapins / Primera
$ws.Cells.Item(314, 1).Value  = 9
$ws.Cells.Item(314, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(314, 3).Value  = "Metropolitana"
$ws.Cells.Item(314, 4).Value  = 44595
$ws.Cells.Item(314, 5).Value  = 13
$ws.Cells.Item(314, 6).Value  = "Fruta"
$ws.Cells.Item(314, 7).Value  = 100103
$ws.Cells.Item(314, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(314, 9).Value  = 100103001
$ws.Cells.Item(314, 10).Value = "Cereza"
$ws.Cells.Item(314, 11).Value = "Lapins"
$ws.Cells.Item(314, 12).Value = "Primera"
$ws.Cells.Item(314, 13).Value = 180
$ws.Cells.Item(314, 14).Value = 10000
$ws.Cells.Item(314, 15).Value = 10000
$ws.Cells.Item(314, 16).Value = 10000
$ws.Cells.Item(314, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(314, 18).Value = "Región del Maule"
$ws.Cells.Item(314, 19).Value = 1000
$ws.Cells.Item(314, 20).Value = 10
